$d = $word.ActiveDocument

# 1) "useEffect( ) will be called once after the component mounted..." paragraph:
#    insert "defined actions/ functions " right after the leading space.
$found1 = $d.Content.Find.Execute(
    " will be called once after the component mounted to the SPA DOM, e.g. when we have navigated to that view. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " defined actions/ functions will be called once after the component mounted to the SPA DOM, e.g. when we have navigated to that view. ",
    2)

# 2) "useEffect( ) will be called again after states or props change..." paragraph:
#    insert "defined actions/ functions " after the leading space, drop the stray
#    "useEffect" before "saying", and replace the trailing "hook is not active then."
#    with "defined  is only called once."
$found2 = $d.Content.Find.Execute(
    " will be called again after states or props change, unless we provide the empty dependency list [ ] to useEffect saying that this useEffect hook is not active then.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " defined actions/ functions will be called again after states or props change, unless we provide the empty dependency list [ ] to saying that this useEffect defined  is only called once.",
    2)

# 3) "/build folder" -> "/dist folder" (keeps the bold run formatting).
$found3 = $d.Content.Find.Execute(
    "/build folder",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "/dist folder",
    2)

Write-Host "Replacement 1 found:" $found1
Write-Host "Replacement 2 found:" $found2
Write-Host "Replacement 3 found:" $found3
